# IWSM_2022_Poster.pptx edit:
#   - "Table 4. "  -> "Table 3. "   (caption above the coefficient table)
#   - "Figure 4. " -> "Figure 3. "  (caption above the cell-type hazard-ratio figure)
#
# Both captions live in slide 1, in text boxes named "TextBox 116" and
# "TextBox 125" respectively. In each case only the number changes; the
# run is therefore split at the digit+period+space so the surrounding
# formatting (Calibri 36pt bold, etc.) is preserved and a fresh run is
# produced for the replaced text, just like typing over a selection in
# the PowerPoint UI would.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Replace-CaptionNumber {
    param(
        [string]$ShapeName,
        [string]$OldNumber,
        [string]$NewNumber
    )

    $shp = $s.Shapes.Item($ShapeName)
    $tr = $shp.TextFrame.TextRange
    $fullText = $tr.Text

    $start = $fullText.IndexOf($OldNumber)
    if ($start -lt 0) {
        throw "Could not find '$OldNumber' in shape '$ShapeName' (text: $fullText)"
    }

    # COM/TextRange.Characters is 1-based.
    $target = $tr.Characters($start + 1, $OldNumber.Length)
    $target.Text = $NewNumber
}

Replace-CaptionNumber "TextBox 116" "4. " "3. "
Replace-CaptionNumber "TextBox 125" "4. " "3. "
